# Add a new transaction row (row 4) to the expense sheet:
#   Time=15:55, Description="Lending a Friend", Category="Lend",
#   Income=0, Expense=400, Balance=-400
# This mirrors rows 2-3: Time/Description/Category/Income/Expense are
# stored as text, Balance is a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "15:55"
$ws.Range("B4").Value = "Lending a Friend"
$ws.Range("C4").Value = "Lend"

# Leading apostrophe forces these numeric-looking entries to stay text,
# matching how Income/Expense are stored as text in the existing rows.
$ws.Range("D4").Value = "'0"
$ws.Range("E4").Value = "'400"

# Balance is a genuine numeric value (Expense - Income sign convention).
$ws.Range("F4").Value = -400
